# feat: add 2022-Q3 data
#
# 1) "总计" (summary) sheet: insert a new row for 2022-Q3 at the top of the
#    data (row 2), pushing the existing quarters down by one row.
# 2) Insert a brand-new "2022-Q3" worksheet right after "总计", carrying the
#    per-fund holding details for the new quarter (two funds).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift the existing 3 data rows (2..4) down to (3..5). Walk back-to-front
# so we never overwrite a row before it has been copied onward.
$total.Range("A4:D4").Copy($total.Range("A5:D5"))
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

# The A column is just a 0-based row index - renumber it after the shift.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3

# Write the new first data row: 2022-Q3, 2 funds held, 0 yi held.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet, right after "总计"
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $total)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Columns B and D:G hold numeric-looking text (fund codes / percentages
# stored as strings in the source data) - force text formatting so the
# values are not silently reinterpreted as numbers.
$q3.Range("B2:B3").NumberFormat = "@"
$q3.Range("D2:G3").NumberFormat = "@"

# Row 2: fund 001375
$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "001375"
$q3.Cells.Item(2, 3).Value = "金元顺安优质精选灵活配置混合C"
$q3.Cells.Item(2, 4).Value = "0.62"
$q3.Cells.Item(2, 5).Value = "65.13"
$q3.Cells.Item(2, 6).Value = "0.72"
$q3.Cells.Item(2, 7).Value = "0.0045"
$q3.Cells.Item(2, 8).Value = 9

# Row 3 (new row): fund 620007
$q3.Range("A2").Copy()
$q3.Range("A3").PasteSpecial(-4122)
$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "620007"
$q3.Cells.Item(3, 3).Value = "金元顺安优质精选灵活配置混合A"
$q3.Cells.Item(3, 4).Value = "0.06"
$q3.Cells.Item(3, 5).Value = "65.13"
$q3.Cells.Item(3, 6).Value = "0.72"
$q3.Cells.Item(3, 7).Value = "0.0004"
$q3.Cells.Item(3, 8).Value = 9

# Restore the original active sheet/tab.
$total.Activate()
